$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so numeric-looking
# strings (e.g. "67.260.79", "0.615") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.260.79'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '3.575.14'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '570.58'
$ws.Range('E5').Value = '  -5.07%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '190.92'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '3.571.28'
$ws.Range('E7').Value = '  -2.69%  '
$ws.Range('D8').Value = '0.615'
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '0.678'
$ws.Range('E10').Value = '  -4.21%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').Value = '56.05'
$ws.Range('E11').Value = '  -3.79%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '0.149'
$ws.Range('E12').Value = '  -2.82%  '
$ws.Range('D13').Value = '0.0000270'
$ws.Range('E13').Value = '  -2.45%  '
$ws.Range('D14').Value = '9.86'
$ws.Range('E14').Value = '  -3.89%  '
$ws.Range('D15').Value = '4.135.03'
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('D16').Value = '3.572.86'
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').Value = '67.095.01'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D19').Value = '12.21'
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('D20').Value = '18.20'
$ws.Range('E20').Value = '  -4.53%  '
$ws.Range('E21').Value = '  -5.60%  '
$ws.Range('D22').Value = '402.35'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '4.15'
$ws.Range('E23').Value = '  -6.99%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').Value = '11.89'
$ws.Range('E24').Value = '  +4.52%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '86.01'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').Value = '2.92'
$ws.Range('E26').Value = '  -2.02%  '
$ws.Range('D27').Value = '12.48'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').Value = '6.10'
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('D29').Value = '3.66'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('D30').Value = '7.79'
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('D31').Value = '8.96'
$ws.Range('E31').Value = '  -4.36%  '
$ws.Range('D32').Value = '31.21'
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').Value = '638.61'
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('D34').Value = '12.12'
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.114'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = '64.00'
$ws.Range('E36').Value = '  -4.78%  '
$ws.Range('D37').Value = '42.42'
$ws.Range('E37').Value = '  -7.45%  '
$ws.Range('D38').Value = '0.405'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').Value = '0.0₃0764'
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('D41').Value = '3.217.60'
$ws.Range('E41').Value = '  +13.35%  '
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('E44').Value = '  +2.13%  '
$ws.Range('D45').Value = '2.71'
$ws.Range('E45').Value = '  +5.05%  '
$ws.Range('D46').Value = '0.0414'
$ws.Range('E46').Value = '  -3.77%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.130'
$ws.Range('E47').Value = '  -4.46%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.08'
$ws.Range('E48').Value = '  -3.73%  '
$ws.Range('D49').Value = '141.89'
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('D50').Value = '8.52'
$ws.Range('E50').Value = '  -5.29%  '
$ws.Range('D51').Value = '2.53'
$ws.Range('E51').Value = '  -4.58%  '

# Restore the original (default) style on the Price column now that the
# text values are set, so no residual formatting differences remain.
$ws.Range("D2:D51").Style = "Normal"

